# se sube para emitir las pólizas de movilidad en QA
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Patente (W2), Motor (X2) and Chasis (Y2) with the new RGM011 values
$ws.Range("W2").Value = "RGM011"
$ws.Range("X2").Value = "1234567RGM011"
$ws.Range("Y2").Value = "1234567RGM011"

# Update SumaAsegurada (U2)
$ws.Range("U2").Value = 1400000

# Move the active selection from F8 to A2
$ws.Range("A2").Select()
